# Generate Report for Handback
# Adds a second handed-back file (ad84fe27-7571-4181-b764-17442849d730) to the
# handback-status report, alongside the existing 62288978-c74f-438f-83ff-b02031d3c663
# entry (renamed from the old bc7ad01a-... UUID), on all three sheets:
#   Overview (sheet1), zh-cn (sheet2), de-de (sheet3).

$wb = $excel.ActiveWorkbook

$oldUuid = "bc7ad01a-d8c4-4593-bace-17fb2811f112"
$uuid1   = "62288978-c74f-438f-83ff-b02031d3c663"
$uuid2   = "ad84fe27-7571-4181-b764-17442849d730"

$hash1zh = "cf34051a15e942c44782187442fcce09126d858c"
$hash1de = "cf34051a15e942c44782187442fcce09126d858c"
$hash2zh = "25dd66467bfead44e6297d6567a5bd2a0349368d"
$hash2de = "25dd66467bfead44e6297d6567a5bd2a0349368d"

$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) as BGR-packed OLE color == FF6495ED

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = 2
    $rng.Font.Color = $hyperlinkColor
}

function Style-AsDate($rng) {
    $rng.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$loOv.ListRows.Add() | Out-Null

$wsOv.Range("A2").Value = "$uuid1.md"
$wsOv.Range("B2").Value = "e2e\$uuid1.md"
$wsOv.Range("G2").Value = "2016-08-13 19:17:53"

$wsOv.Range("A3").Value = "$uuid2.md"
$wsOv.Range("B3").Value = "e2e\$uuid2.md"
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("E3").Value = "Handed back: in sync with en-US"
$wsOv.Range("F3").Value = "Handed back: in sync with en-US"
$wsOv.Range("G3").Value = "2016-08-13 19:17:53"

Style-AsHyperlink $wsOv.Range("B3")
Style-AsDate $wsOv.Range("G3")

$wsOv.Hyperlinks.Delete()
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/203d4f62eaec768c675b42d1c701148cc6893d7a/e2e/$uuid1.md", "", "", "e2e\$uuid1.md") | Out-Null
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/203d4f62eaec768c675b42d1c701148cc6893d7a/e2e/$uuid2.md", "", "", "e2e\$uuid2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A2").Value = "$uuid1.md"
$wsZh.Range("B2").Value = ".md"
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "e2e"
$wsZh.Range("E2").Value = "ht"
$wsZh.Range("F2").Value = "False"
$wsZh.Range("G2").Value = "$uuid1.$hash1zh.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-13 19:17:46"
$wsZh.Range("I2").Value = "$uuid1.md"
$wsZh.Range("J2").Value = "$uuid1.$hash1zh.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-13 19:18:29"
$wsZh.Range("L2").Value = ""
$wsZh.Range("M2").Value = "True"
$wsZh.Range("N2").Value = ""
$wsZh.Range("O2").Value = "False"
$wsZh.Range("P2").Value = ""

$wsZh.Range("A3").Value = "$uuid2.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$uuid2.$hash2zh.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-13 19:17:46"
$wsZh.Range("I3").Value = "$uuid2.md"
$wsZh.Range("J3").Value = "$uuid2.$hash2zh.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-13 19:18:29"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

Style-AsHyperlink $wsZh.Range("A3")
Style-AsHyperlink $wsZh.Range("I3")
Style-AsDate $wsZh.Range("H2")
Style-AsDate $wsZh.Range("K2")
Style-AsDate $wsZh.Range("H3")
Style-AsDate $wsZh.Range("K3")

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/203d4f62eaec768c675b42d1c701148cc6893d7a/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/78ecf3e47a59c97e80f0fb2976d05c732dc935ff/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/203d4f62eaec768c675b42d1c701148cc6893d7a/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/78ecf3e47a59c97e80f0fb2976d05c732dc935ff/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A2").Value = "$uuid1.md"
$wsDe.Range("B2").Value = ".md"
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "e2e"
$wsDe.Range("E2").Value = "ht"
$wsDe.Range("F2").Value = "False"
$wsDe.Range("G2").Value = "$uuid1.$hash1de.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-13 19:17:53"
$wsDe.Range("I2").Value = "$uuid1.md"
$wsDe.Range("J2").Value = "$uuid1.$hash1de.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-13 19:18:39"
$wsDe.Range("L2").Value = ""
$wsDe.Range("M2").Value = "True"
$wsDe.Range("N2").Value = ""
$wsDe.Range("O2").Value = "False"
$wsDe.Range("P2").Value = ""

$wsDe.Range("A3").Value = "$uuid2.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$uuid2.$hash2de.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-13 19:17:53"
$wsDe.Range("I3").Value = "$uuid2.md"
$wsDe.Range("J3").Value = "$uuid2.$hash2de.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-13 19:18:39"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

Style-AsHyperlink $wsDe.Range("A3")
Style-AsHyperlink $wsDe.Range("I3")
Style-AsDate $wsDe.Range("H2")
Style-AsDate $wsDe.Range("K2")
Style-AsDate $wsDe.Range("H3")
Style-AsDate $wsDe.Range("K3")

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/203d4f62eaec768c675b42d1c701148cc6893d7a/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8ab3d00caf3d01b40a8c262eb1be61ee025db10c/e2e/$uuid1.md", "", "", "$uuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/203d4f62eaec768c675b42d1c701148cc6893d7a/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8ab3d00caf3d01b40a8c262eb1be61ee025db10c/e2e/$uuid2.md", "", "", "$uuid2.md") | Out-Null

Write-Output "done"
